# Generate Report for Handoff
# Adds a new "26e9d434-ed36-44d4-a2a1-c4404329daa4.md" row above the
# existing "f2adc214-cd29-47f6-b049-ca1d6554e105.md" row on every sheet
# (Overview, zh-cn, de-de), pushing the old row down, and keeps the
# ListObjects / dimensions / hyperlinks in sync.

$wb = $excel.ActiveWorkbook

$baseUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/4a4933ec60b2b233f4be10fd5f33dbb9ea0b36f0/e2e/"
$newFile = "26e9d434-ed36-44d4-a2a1-c4404329daa4.md"
$oldFile = "f2adc214-cd29-47f6-b049-ca1d6554e105.md"

# ----------------------------------------------------------------------
# Sheet 1: Overview
# ----------------------------------------------------------------------
$wsOv = $wb.Worksheets.Item("Overview")

$wsOv.Rows.Item(2).Insert()
$wsOv.Range("B2").Hyperlinks.Delete()

$wsOv.Cells.Item(2,1).Value = $newFile
$wsOv.Cells.Item(2,2).Value = "e2e\" + $newFile
$wsOv.Cells.Item(2,3).Value = ".md"
$wsOv.Cells.Item(2,4).NumberFormat = "General"
$wsOv.Cells.Item(2,4).Value = ""
$wsOv.Cells.Item(2,5).Value = "Ready for handoff"
$wsOv.Cells.Item(2,6).Value = "Ready for handoff"
$wsOv.Cells.Item(2,7).Value = "2016-09-05 18:45:18"
$wsOv.Cells.Item(2,7).NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsOv.Hyperlinks.Add($wsOv.Range("B2"), ($baseUrl + $newFile), "", "", ("e2e\" + $newFile)) | Out-Null
$wsOv.Cells.Item(2,2).Style = "HyperLink"

$wsOv.Hyperlinks.Add($wsOv.Range("B3"), ($baseUrl + $oldFile), "", "", ("e2e\" + $oldFile)) | Out-Null
$wsOv.Cells.Item(3,2).Style = "HyperLink"

$loOv = $wsOv.ListObjects.Item(1)
$loOv.Resize($wsOv.Range("A1:G3"))

# ----------------------------------------------------------------------
# Sheet 2: zh-cn
# ----------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Rows.Item(2).Insert()
$wsZh.Range("A2").Hyperlinks.Delete()

$wsZh.Cells.Item(2,1).Value = $newFile
$wsZh.Cells.Item(2,2).Value = ".md"
$wsZh.Cells.Item(2,3).Value = "Ready for handoff"
$wsZh.Cells.Item(2,4).Value = "e2e"
$wsZh.Cells.Item(2,5).Value = "ht"
$wsZh.Cells.Item(2,6).Value = "'False"
$wsZh.Cells.Item(2,7).Value = "26e9d434-ed36-44d4-a2a1-c4404329daa4.aa7e8ced1d1dc7747f813b64d93d8b23b6343f3f.zh-cn.xlf"
$wsZh.Cells.Item(2,8).Value = "2016-09-05 18:45:12"
$wsZh.Cells.Item(2,8).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Cells.Item(2,9).NumberFormat = "General"
$wsZh.Cells.Item(2,9).Value = ""
$wsZh.Cells.Item(2,10).NumberFormat = "General"
$wsZh.Cells.Item(2,10).Value = ""
$wsZh.Cells.Item(2,11).Value = "0001-01-01 00:00:00"
$wsZh.Cells.Item(2,11).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Cells.Item(2,12).NumberFormat = "General"
$wsZh.Cells.Item(2,12).Value = ""
$wsZh.Cells.Item(2,13).Value = "'True"
$wsZh.Cells.Item(2,14).NumberFormat = "General"
$wsZh.Cells.Item(2,14).Value = ""
$wsZh.Cells.Item(2,15).Value = "'False"
$wsZh.Cells.Item(2,16).NumberFormat = "General"
$wsZh.Cells.Item(2,16).Value = ""

$wsZh.Hyperlinks.Add($wsZh.Range("A2"), ($baseUrl + $newFile), "", "", $newFile) | Out-Null
$wsZh.Cells.Item(2,1).Style = "HyperLink"

$wsZh.Hyperlinks.Add($wsZh.Range("A3"), ($baseUrl + $oldFile), "", "", $oldFile) | Out-Null
$wsZh.Cells.Item(3,1).Style = "HyperLink"

$loZh = $wsZh.ListObjects.Item(1)
$loZh.Resize($wsZh.Range("A1:P3"))

# ----------------------------------------------------------------------
# Sheet 3: de-de
# ----------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Rows.Item(2).Insert()
$wsDe.Range("A2").Hyperlinks.Delete()

$wsDe.Cells.Item(2,1).Value = $newFile
$wsDe.Cells.Item(2,2).Value = ".md"
$wsDe.Cells.Item(2,3).Value = "Ready for handoff"
$wsDe.Cells.Item(2,4).Value = "e2e"
$wsDe.Cells.Item(2,5).Value = "ht"
$wsDe.Cells.Item(2,6).Value = "'False"
$wsDe.Cells.Item(2,7).Value = "26e9d434-ed36-44d4-a2a1-c4404329daa4.aa7e8ced1d1dc7747f813b64d93d8b23b6343f3f.de-de.xlf"
$wsDe.Cells.Item(2,8).Value = "2016-09-05 18:45:18"
$wsDe.Cells.Item(2,8).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Cells.Item(2,9).NumberFormat = "General"
$wsDe.Cells.Item(2,9).Value = ""
$wsDe.Cells.Item(2,10).NumberFormat = "General"
$wsDe.Cells.Item(2,10).Value = ""
$wsDe.Cells.Item(2,11).Value = "0001-01-01 00:00:00"
$wsDe.Cells.Item(2,11).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Cells.Item(2,12).NumberFormat = "General"
$wsDe.Cells.Item(2,12).Value = ""
$wsDe.Cells.Item(2,13).Value = "'True"
$wsDe.Cells.Item(2,14).NumberFormat = "General"
$wsDe.Cells.Item(2,14).Value = ""
$wsDe.Cells.Item(2,15).Value = "'False"
$wsDe.Cells.Item(2,16).NumberFormat = "General"
$wsDe.Cells.Item(2,16).Value = ""

$wsDe.Hyperlinks.Add($wsDe.Range("A2"), ($baseUrl + $newFile), "", "", $newFile) | Out-Null
$wsDe.Cells.Item(2,1).Style = "HyperLink"

$wsDe.Hyperlinks.Add($wsDe.Range("A3"), ($baseUrl + $oldFile), "", "", $oldFile) | Out-Null
$wsDe.Cells.Item(3,1).Style = "HyperLink"

$loDe = $wsDe.ListObjects.Item(1)
$loDe.Resize($wsDe.Range("A1:P3"))

Write-Output "done"
